$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Periodo Mora" / "Valor Mora" data rows ---------------------
# Row 16 and Row 17 swap their "Periodo Mora" period labels (2501 <-> 2412)
# together with a corresponding swap of the "Valor Mora" amounts, so the
# newest period (2412) now appears first.
$ws.Range("E16").Value2 = "2412"
$ws.Range("F16").Value2 = 143600

$ws.Range("E17").Value2 = "2501"
$ws.Range("F17").Value2 = 14360

# --- Close off the label boxes with a right border --------------------------
# The logo box (merged B2:C5) and the RAZON SOCIAL / NIT / VALOR MORA label
# cells (merged B7:D7, B9:D9, B11:D11) previously had no right-hand border;
# add a thin right border so each box is fully enclosed on all four sides,
# matching the left/top/bottom thin borders already present.

# B2 already carries the fully-enclosed box style that C2:C5 should match,
# so copy its format across rather than building a brand-new style.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("C2:C5").PasteSpecial(-4122) | Out-Null

# D7 has no existing "donor" cell with a matching right border, so add the
# border directly, then propagate the resulting format to D9 and D11 so all
# three share a single new style instead of each creating its own.
$ws.Range("D7").Borders.Item(10).LineStyle = 1
$ws.Range("D7").Borders.Item(10).Weight = 2

$ws.Range("D7").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4122) | Out-Null
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
